$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.066.16'
$ws.Range("E2").Value = '  -1.93%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.789.58'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.97'
$ws.Range("E5").Value = '  -0.89%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.47'
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("E10").Value = '  -0.94%  '

$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.047.54'
$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.796.72'
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.83'
$ws.Range("E14").Value = '  -1.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.626'
$ws.Range("E15").Value = '  -1.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.078.59'
$ws.Range("E16").Value = '  -1.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.16'
$ws.Range("E17").Value = '  -3.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.09'
$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.30'
$ws.Range("E19").Value = '  -3.94%  '

$ws.Range("E20").Value = '  -3.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.77'
$ws.Range("E22").Value = '  +0.24%  '

$ws.Range("E23").Value = '  -3.04%  '

$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.39'
$ws.Range("E25").Value = '  -1.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.37'
$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.06'
$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("E28").Value = '  -1.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("E30").Value = '  -1.74%  '

$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("E32").Value = '  -3.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.49'
$ws.Range("E33").Value = '  -3.45%  '

$ws.Range("E34").Value = '  -2.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.396.62'
$ws.Range("E35").Value = '  -2.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.639'
$ws.Range("E36").Value = '  +0.76%  '

$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("E38").Value = '  -3.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '79.59'
$ws.Range("E39").Value = '  -6.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.922'
$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("E41").Value = '  +1.01%  '

$ws.Range("E42").Value = '  -3.02%  '

$ws.Range("E43").Value = '  +1.89%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0496'
$ws.Range("E44").Value = '  +0.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.90'
$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.40'
$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("E47").Value = '  -1.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.946.86'
$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.03'
$ws.Range("E49").Value = '  +0.66%  '

$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0128'
$ws.Range("E51").Value = '  +2.35%  '
